# Updates FFXIV Leve market-price snapshot values (currentAveragePrice*, Leve cost, and
# profit columns H/I/J/K/L/M/N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to
# reflect the latest scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 294.36923
$ws.Range("J17").Value = 294.36923
$ws.Range("L17").Value = 883.10769
$ws.Range("N17").Value = -1219.10769

# Row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 4142.857
$ws.Range("I74").Value = 4000
$ws.Range("K74").Value = 4000
$ws.Range("M74").Value = -3064

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 90912070
$ws.Range("I76").Value = 333335680
$ws.Range("J76").Value = 3225
$ws.Range("K76").Value = 333335680
$ws.Range("L76").Value = 3225
$ws.Range("M76").Value = -333335365
$ws.Range("N76").Value = -3855

# Row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 4142.857
$ws.Range("I77").Value = 4000
$ws.Range("K77").Value = 20000
$ws.Range("M77").Value = -15320

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 90912070
$ws.Range("I79").Value = 333335680
$ws.Range("J79").Value = 3225
$ws.Range("K79").Value = 333335680
$ws.Range("L79").Value = 3225
$ws.Range("M79").Value = -333334588
$ws.Range("N79").Value = -5409

# Row 118: Crafty Concoctions
$ws.Range("H118").Value = 3327.3333
$ws.Range("J118").Value = 5639.6
$ws.Range("L118").Value = 16918.8
$ws.Range("N118").Value = -20232.8

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 3368134.2
$ws.Range("I132").Value = 1187.7742
$ws.Range("J132").Value = 55555810
$ws.Range("K132").Value = 3563.3226
$ws.Range("L132").Value = 166667430
$ws.Range("M132").Value = -1033.3226
$ws.Range("N132").Value = -166672490

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 75049820
$ws.Range("I137").Value = 41667570
$ws.Range("J137").Value = 89356504
$ws.Range("K137").Value = 125002710
$ws.Range("L137").Value = 268069512
$ws.Range("M137").Value = -125000160
$ws.Range("N137").Value = -268074612

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 7189653
$ws.Range("I32").Value = 1816304.6
$ws.Range("K32").Value = 1816304.6
$ws.Range("M32").Value = -1816017.6

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 358264.56
$ws.Range("I45").Value = 667627.6
$ws.Range("J45").Value = 1307.2307
$ws.Range("K45").Value = 667627.6
$ws.Range("L45").Value = 1307.2307
$ws.Range("M45").Value = -667250.6
$ws.Range("N45").Value = -2061.2307

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 1988.7778
$ws.Range("I63").Value = 1983.3334
$ws.Range("J63").Value = 1999.6666
$ws.Range("K63").Value = 1983.3334
$ws.Range("L63").Value = 1999.6666
$ws.Range("M63").Value = -1297.3334
$ws.Range("N63").Value = -3371.6666

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 1988.7778
$ws.Range("I66").Value = 1983.3334
$ws.Range("J66").Value = 1999.6666
$ws.Range("K66").Value = 9916.666999999999
$ws.Range("L66").Value = 9998.333000000001
$ws.Range("M66").Value = -6484.666999999999
$ws.Range("N66").Value = -16862.333

# Row 88: The Mast Chance
$ws.Range("H88").Value = 6556.25
$ws.Range("J88").Value = 6556.25
$ws.Range("L88").Value = 6556.25
$ws.Range("N88").Value = -7368.25

# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 6556.25
$ws.Range("J91").Value = 6556.25
$ws.Range("L91").Value = 6556.25
$ws.Range("N91").Value = -9364.25

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1786.6552
$ws.Range("I102").Value = 1781.2693
$ws.Range("J102").Value = 1833.3334
$ws.Range("K102").Value = 1781.2693
$ws.Range("L102").Value = 1833.3334
$ws.Range("M102").Value = -159.2692999999999
$ws.Range("N102").Value = -5077.3334

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1081.4286
$ws.Range("I110").Value = 1137.4
$ws.Range("J110").Value = 941.5
$ws.Range("K110").Value = 1137.4
$ws.Range("L110").Value = 941.5
$ws.Range("M110").Value = 907.5999999999999
$ws.Range("N110").Value = -5031.5

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 9974962
$ws.Range("I132").Value = 9527362
$ws.Range("J132").Value = 13891464
$ws.Range("K132").Value = 28582086
$ws.Range("L132").Value = 41674392
$ws.Range("M132").Value = -28579556
$ws.Range("N132").Value = -41679452

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 1965.64
$ws.Range("I86").Value = 1972.9166
$ws.Range("J86").Value = 1791
$ws.Range("K86").Value = 1972.9166
$ws.Range("L86").Value = 1791
$ws.Range("M86").Value = -849.9166
$ws.Range("N86").Value = -4037

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1965.64
$ws.Range("I89").Value = 1972.9166
$ws.Range("J89").Value = 1791
$ws.Range("K89").Value = 9864.583000000001
$ws.Range("L89").Value = 8955
$ws.Range("M89").Value = -4248.583000000001
$ws.Range("N89").Value = -20187

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2200
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -6494

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 7339793
$ws.Range("I134").Value = 7937529
$ws.Range("J134").Value = 3574056.5
$ws.Range("K134").Value = 23812587
$ws.Range("L134").Value = 10722169.5
$ws.Range("M134").Value = -23810052
$ws.Range("N134").Value = -10727239.5

$ws = $wb.Worksheets.Item("CRP")
# Row 51: Greenstone for Greenhorns
$ws.Range("H51").Value = 8324.25
$ws.Range("J51").Value = 8324.25
$ws.Range("L51").Value = 8324.25
$ws.Range("N51").Value = -9796.25

# Row 61: Incant Now, Think Later
$ws.Range("H61").Value = 8324.25
$ws.Range("J61").Value = 8324.25
$ws.Range("L61").Value = 8324.25
$ws.Range("N61").Value = -9020.25

# Row 68: Do You Even String Bow
$ws.Range("H68").Value = 29721.25
$ws.Range("J68").Value = 30295
$ws.Range("L68").Value = 30295
$ws.Range("N68").Value = -31793

# Row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value = 29721.25
$ws.Range("J71").Value = 30295
$ws.Range("L71").Value = 90885
$ws.Range("N71").Value = -98373

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 4978.5454
$ws.Range("I122").Value = 8561.75
$ws.Range("J122").Value = 678.7
$ws.Range("K122").Value = 25685.25
$ws.Range("L122").Value = 2036.1
$ws.Range("M122").Value = -23235.25
$ws.Range("N122").Value = -6936.1

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1564911.2
$ws.Range("I132").Value = 2274375
$ws.Range("J132").Value = 4091
$ws.Range("K132").Value = 6823125
$ws.Range("L132").Value = 12273
$ws.Range("M132").Value = -6820595
$ws.Range("N132").Value = -17333

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1147848.1
$ws.Range("I134").Value = 5088.1924
$ws.Range("J134").Value = 4449154.5
$ws.Range("K134").Value = 15264.5772
$ws.Range("L134").Value = 13347463.5
$ws.Range("M134").Value = -12729.5772
$ws.Range("N134").Value = -13352533.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 5803956.5
$ws.Range("I5").Value = 6993898.5
$ws.Range("K5").Value = 20981695.5
$ws.Range("M5").Value = -20981583.5

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 1128.4642
$ws.Range("J113").Value = 1753.9642
$ws.Range("L113").Value = 5261.892599999999
$ws.Range("N113").Value = -9601.892599999999

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 924.7838
$ws.Range("I131").Value = 535
$ws.Range("J131").Value = 972.0303
$ws.Range("K131").Value = 1605
$ws.Range("L131").Value = 2916.0909
$ws.Range("M131").Value = 3435
$ws.Range("N131").Value = -12996.0909

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 5803956.5
$ws.Range("I135").Value = 6993898.5
$ws.Range("K135").Value = 62945086.5
$ws.Range("M135").Value = -62942551.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 2718416.5
$ws.Range("I70").Value = 1689283
$ws.Range("J70").Value = 4335626.5
$ws.Range("K70").Value = 1689283
$ws.Range("L70").Value = 4335626.5
$ws.Range("M70").Value = -1689013
$ws.Range("N70").Value = -4336166.5

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 2718416.5
$ws.Range("I73").Value = 1689283
$ws.Range("J73").Value = 4335626.5
$ws.Range("K73").Value = 1689283
$ws.Range("L73").Value = 4335626.5
$ws.Range("M73").Value = -1688347
$ws.Range("N73").Value = -4337498.5

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 16077.177
$ws.Range("I80").Value = 30000
$ws.Range("J80").Value = 15207
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 15207
$ws.Range("M80").Value = -29002
$ws.Range("N80").Value = -17203

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 16077.177
$ws.Range("I83").Value = 30000
$ws.Range("J83").Value = 15207
$ws.Range("K83").Value = 150000
$ws.Range("L83").Value = 76035
$ws.Range("M83").Value = -145008
$ws.Range("N83").Value = -86019

# Row 132: On Board for Lar
$ws.Range("H132").Value = 30778558
$ws.Range("I132").Value = 41272010
$ws.Range("J132").Value = 18186416
$ws.Range("K132").Value = 123816030
$ws.Range("L132").Value = 54559248
$ws.Range("M132").Value = -123813500
$ws.Range("N132").Value = -54564308

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 1139.4445
$ws.Range("I7").Value = 942.1667
$ws.Range("J7").Value = 1534
$ws.Range("K7").Value = 942.1667
$ws.Range("L7").Value = 1534
$ws.Range("M7").Value = -830.1667
$ws.Range("N7").Value = -1758

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 603.9167
$ws.Range("I46").Value = 777
$ws.Range("J46").Value = 430.83334
$ws.Range("K46").Value = 777
$ws.Range("L46").Value = 430.83334
$ws.Range("M46").Value = -589
$ws.Range("N46").Value = -806.83334

# Row 126: Battered Books
$ws.Range("H126").Value = 1139.4445
$ws.Range("I126").Value = 942.1667
$ws.Range("J126").Value = 1534
$ws.Range("K126").Value = 2826.5001
$ws.Range("L126").Value = 4602
$ws.Range("M126").Value = -356.5001000000002
$ws.Range("N126").Value = -9542

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 1338104
$ws.Range("I136").Value = 1839348.2
$ws.Range("K136").Value = 5518044.6
$ws.Range("M136").Value = -5515494.6

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1622444.8
$ws.Range("I132").Value = 1234274.2
$ws.Range("K132").Value = 3702822.6
$ws.Range("M132").Value = -3700292.6
